# The deck currently uses the "Integral" theme (custom green palette) for
# its slide master / main presentation theme, while a second, unused
# "Office Theme" (default blue/orange palette) only backs the Notes
# Master. The edit swaps the two: the main presentation theme becomes the
# plain default "Office Theme" palette, while the (otherwise unused)
# secondary theme slot ends up holding the "Integral" palette.
#
# Practically, the PowerPoint object model lets us repaint the live
# theme's 12-slot colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) via ThemeColorScheme on a slide. We push the default Office
# theme colours into that scheme so the active theme part matches the
# "Office Theme" palette exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToOle([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2,
# 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Colors(1).RGB  = ToOle 0x00 0x00 0x00   # dk1
$tcs.Colors(2).RGB  = ToOle 0xFF 0xFF 0xFF   # lt1
$tcs.Colors(3).RGB  = ToOle 0x44 0x54 0x6A   # dk2
$tcs.Colors(4).RGB  = ToOle 0xE7 0xE6 0xE6   # lt2
$tcs.Colors(5).RGB  = ToOle 0x5B 0x9B 0xD5   # accent1
$tcs.Colors(6).RGB  = ToOle 0xED 0x7D 0x31   # accent2
$tcs.Colors(7).RGB  = ToOle 0xA5 0xA5 0xA5   # accent3
$tcs.Colors(8).RGB  = ToOle 0xFF 0xC0 0x00   # accent4
$tcs.Colors(9).RGB  = ToOle 0x44 0x72 0xC4   # accent5
$tcs.Colors(10).RGB = ToOle 0x70 0xAD 0x47   # accent6
$tcs.Colors(11).RGB = ToOle 0x05 0x63 0xC1   # hlink
$tcs.Colors(12).RGB = ToOle 0x95 0x4F 0x72   # folHlink
